# Auto-generated script: applies the 2023-10-31 data update
# For each affected worksheet/row, update column J (2023 total) to its new value.
$wb = $excel.ActiveWorkbook


# --- Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 10).Value = 6412  # Aggravated Assault: 6399 -> 6412
$ws.Cells.Item(3, 10).Value = 6813  # Aggravated Battery: 6795 -> 6813
$ws.Cells.Item(4, 10).Value = 1472  # Criminal Sexual Assault: 1466 -> 1472
$ws.Cells.Item(5, 10).Value = 525  # Homicide: 522 -> 525
$ws.Cells.Item(6, 10).Value = 9031  # Robbery: 9006 -> 9031
$ws.Cells.Item(7, 10).Value = 24253  # Total: 24188 -> 24253

# --- Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(3, 10).Value = 460  # Aggravated Battery: 458 -> 460
$ws.Cells.Item(6, 10).Value = 539  # Robbery: 536 -> 539
$ws.Cells.Item(7, 10).Value = 1533  # Total: 1528 -> 1533

# --- South Chicago ---
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 10).Value = 142  # Aggravated Assault: 141 -> 142
$ws.Cells.Item(3, 10).Value = 180  # Aggravated Battery: 179 -> 180
$ws.Cells.Item(6, 10).Value = 128  # Robbery: 127 -> 128
$ws.Cells.Item(7, 10).Value = 481  # Total: 478 -> 481

# --- Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(4, 10).Value = 46  # Criminal Sexual Assault: 45 -> 46
$ws.Cells.Item(7, 10).Value = 1090  # Total: 1089 -> 1090

# --- West Pullman ---
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(3, 10).Value = 127  # Aggravated Battery: 126 -> 127
$ws.Cells.Item(7, 10).Value = 350  # Total: 349 -> 350

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(3, 10).Value = 253  # Aggravated Battery: 252 -> 253
$ws.Cells.Item(5, 10).Value = 29  # Homicide: 28 -> 29
$ws.Cells.Item(6, 10).Value = 220  # Robbery: 219 -> 220
$ws.Cells.Item(7, 10).Value = 745  # Total: 742 -> 745

# --- Woodlawn ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(6, 10).Value = 95  # Robbery: 94 -> 95
$ws.Cells.Item(7, 10).Value = 372  # Total: 371 -> 372

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(2, 10).Value = 193  # Albany Park: 192 -> 193
$ws.Cells.Item(6, 10).Value = 183  # Ashburn: 182 -> 183
$ws.Cells.Item(7, 10).Value = 706  # Auburn Gresham: 704 -> 706
$ws.Cells.Item(8, 10).Value = 1533  # Austin: 1528 -> 1533
$ws.Cells.Item(9, 10).Value = 125  # Avalon Park: 124 -> 125
$ws.Cells.Item(11, 10).Value = 409  # Belmont Cragin: 408 -> 409
$ws.Cells.Item(14, 10).Value = 128  # Bridgeport: 126 -> 128
$ws.Cells.Item(18, 10).Value = 208  # Calumet Heights: 207 -> 208
$ws.Cells.Item(19, 10).Value = 709  # Chatham: 708 -> 709
$ws.Cells.Item(20, 10).Value = 505  # Chicago Lawn: 504 -> 505
$ws.Cells.Item(25, 10).Value = 119  # East Side: 118 -> 119
$ws.Cells.Item(29, 10).Value = 1323  # Englewood: 1322 -> 1323
$ws.Cells.Item(31, 10).Value = 229  # Gage Park: 228 -> 229
$ws.Cells.Item(33, 10).Value = 1090  # Garfield Park: 1089 -> 1090
$ws.Cells.Item(36, 10).Value = 330  # Grand Boulevard: 329 -> 330
$ws.Cells.Item(37, 10).Value = 745  # Grand Crossing: 742 -> 745
$ws.Cells.Item(42, 10).Value = 1045  # Humboldt Park: 1043 -> 1045
$ws.Cells.Item(49, 10).Value = 154  # Lincoln Park: 153 -> 154
$ws.Cells.Item(50, 10).Value = 147  # Lincoln Square: 146 -> 147
$ws.Cells.Item(51, 10).Value = 299  # Little Italy, UIC: 297 -> 299
$ws.Cells.Item(52, 10).Value = 613  # Little Village: 612 -> 613
$ws.Cells.Item(57, 10).Value = 107  # Mckinley Park: 106 -> 107
$ws.Cells.Item(63, 10).Value = 83  # NO NEIGHBORHOOD DATA: 80 -> 83
$ws.Cells.Item(64, 10).Value = 157  # Near South Side: 156 -> 157
$ws.Cells.Item(67, 10).Value = 913  # North Lawndale: 908 -> 913
$ws.Cells.Item(73, 10).Value = 233  # Portage Park: 232 -> 233
$ws.Cells.Item(76, 10).Value = 362  # River North: 360 -> 362
$ws.Cells.Item(78, 10).Value = 287  # Rogers Park: 285 -> 287
$ws.Cells.Item(79, 10).Value = 684  # Roseland: 682 -> 684
$ws.Cells.Item(83, 10).Value = 481  # South Chicago: 478 -> 481
$ws.Cells.Item(84, 10).Value = 201  # South Deering: 200 -> 201
$ws.Cells.Item(85, 10).Value = 1010  # South Shore: 1004 -> 1010
$ws.Cells.Item(89, 10).Value = 315  # Uptown: 314 -> 315
$ws.Cells.Item(91, 10).Value = 277  # Washington Park: 274 -> 277
$ws.Cells.Item(95, 10).Value = 350  # West Pullman: 349 -> 350
$ws.Cells.Item(98, 10).Value = 180  # Wicker Park: 179 -> 180
$ws.Cells.Item(99, 10).Value = 372  # Woodlawn: 371 -> 372
$ws.Cells.Item(101, 10).Value = 24253  # Total: 24188 -> 24253

# --- Gage Park ---
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(2, 10).Value = 83  # Aggravated Assault: 82 -> 83
$ws.Cells.Item(7, 10).Value = 229  # Total: 228 -> 229

# --- North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 10).Value = 230  # Aggravated Assault: 227 -> 230
$ws.Cells.Item(3, 10).Value = 340  # Aggravated Battery: 339 -> 340
$ws.Cells.Item(6, 10).Value = 252  # Robbery: 251 -> 252
$ws.Cells.Item(7, 10).Value = 913  # Total: 908 -> 913

# --- South Deering ---
$ws = $wb.Worksheets.Item("South Deering")
$ws.Cells.Item(3, 10).Value = 66  # Aggravated Battery: 65 -> 66
$ws.Cells.Item(7, 10).Value = 201  # Total: 200 -> 201

# --- Lincoln Park ---
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Cells.Item(2, 10).Value = 27  # Aggravated Assault: 26 -> 27
$ws.Cells.Item(7, 10).Value = 154  # Total: 153 -> 154

# --- Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(6, 10).Value = 333  # Robbery: 332 -> 333
$ws.Cells.Item(7, 10).Value = 1323  # Total: 1322 -> 1323

# --- Chatham ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(3, 10).Value = 206  # Aggravated Battery: 205 -> 206
$ws.Cells.Item(7, 10).Value = 709  # Total: 708 -> 709

# --- River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(3, 10).Value = 73  # Aggravated Battery: 72 -> 73
$ws.Cells.Item(6, 10).Value = 199  # Robbery: 198 -> 199
$ws.Cells.Item(7, 10).Value = 362  # Total: 360 -> 362

# --- Bridgeport ---
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Cells.Item(2, 10).Value = 43  # Aggravated Assault: 42 -> 43
$ws.Cells.Item(6, 10).Value = 50  # Robbery: 49 -> 50
$ws.Cells.Item(7, 10).Value = 128  # Total: 126 -> 128

# --- Ashburn ---
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Cells.Item(4, 10).Value = 11  # Criminal Sexual Assault: 10 -> 11
$ws.Cells.Item(7, 10).Value = 183  # Total: 182 -> 183

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(3, 10).Value = 202  # Aggravated Battery: 201 -> 202
$ws.Cells.Item(5, 10).Value = 20  # Homicide: 19 -> 20
$ws.Cells.Item(7, 10).Value = 1045  # Total: 1043 -> 1045

# --- Rogers Park ---
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(2, 10).Value = 77  # Aggravated Assault: 76 -> 77
$ws.Cells.Item(6, 10).Value = 85  # Robbery: 84 -> 85
$ws.Cells.Item(7, 10).Value = 287  # Total: 285 -> 287

# --- Washington Park ---
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(3, 10).Value = 115  # Aggravated Battery: 113 -> 115
$ws.Cells.Item(6, 10).Value = 69  # Robbery: 68 -> 69
$ws.Cells.Item(7, 10).Value = 277  # Total: 274 -> 277

# --- Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 10).Value = 190  # Aggravated Assault: 189 -> 190
$ws.Cells.Item(3, 10).Value = 230  # Aggravated Battery: 229 -> 230
$ws.Cells.Item(7, 10).Value = 684  # Total: 682 -> 684

# --- Near South Side ---
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Cells.Item(2, 10).Value = 43  # Aggravated Assault: 42 -> 43
$ws.Cells.Item(7, 10).Value = 157  # Total: 156 -> 157

# --- Chicago Lawn ---
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(6, 10).Value = 137  # Robbery: 136 -> 137
$ws.Cells.Item(7, 10).Value = 505  # Total: 504 -> 505

# --- Calumet Heights ---
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(5, 10).Value = 3  # Homicide: 2 -> 3
$ws.Cells.Item(7, 10).Value = 208  # Total: 207 -> 208

# --- Grand Boulevard ---
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(6, 10).Value = 100  # Robbery: 99 -> 100
$ws.Cells.Item(7, 10).Value = 330  # Total: 329 -> 330

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(6, 10).Value = 227  # Robbery: 225 -> 227
$ws.Cells.Item(7, 10).Value = 706  # Total: 704 -> 706

# --- East Side ---
$ws = $wb.Worksheets.Item("East Side")
$ws.Cells.Item(3, 10).Value = 35  # Aggravated Battery: 34 -> 35
$ws.Cells.Item(7, 10).Value = 119  # Total: 118 -> 119

# --- Wicker Park ---
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Cells.Item(2, 10).Value = 33  # Aggravated Assault: 32 -> 33
$ws.Cells.Item(7, 10).Value = 180  # Total: 179 -> 180

# --- Lincoln Square ---
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Cells.Item(6, 10).Value = 50  # Robbery: 49 -> 50
$ws.Cells.Item(7, 10).Value = 147  # Total: 146 -> 147

# --- Belmont Cragin ---
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(6, 10).Value = 184  # Robbery: 183 -> 184
$ws.Cells.Item(7, 10).Value = 409  # Total: 408 -> 409

# --- Avalon Park ---
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(6, 10).Value = 43  # Robbery: 42 -> 43
$ws.Cells.Item(7, 10).Value = 125  # Total: 124 -> 125

# --- Portage Park ---
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(3, 10).Value = 60  # Aggravated Battery: 59 -> 60
$ws.Cells.Item(7, 10).Value = 233  # Total: 232 -> 233

# --- Albany Park ---
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Cells.Item(3, 10).Value = 48  # Aggravated Battery: 47 -> 48
$ws.Cells.Item(7, 10).Value = 193  # Total: 192 -> 193

# --- Uptown ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(4, 10).Value = 32  # Criminal Sexual Assault: 31 -> 32
$ws.Cells.Item(7, 10).Value = 315  # Total: 314 -> 315

# --- Little Italy, UIC ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(2, 10).Value = 67  # Aggravated Assault: 66 -> 67
$ws.Cells.Item(4, 10).Value = 26  # Criminal Sexual Assault: 25 -> 26
$ws.Cells.Item(7, 10).Value = 299  # Total: 297 -> 299

# --- Mckinley Park ---
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Cells.Item(6, 10).Value = 48  # Robbery: 47 -> 48
$ws.Cells.Item(7, 10).Value = 107  # Total: 106 -> 107

# --- South Shore ---
$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(2, 10).Value = 268  # Aggravated Assault: 266 -> 268
$ws.Cells.Item(3, 10).Value = 358  # Aggravated Battery: 357 -> 358
$ws.Cells.Item(4, 10).Value = 67  # Criminal Sexual Assault: 66 -> 67
$ws.Cells.Item(6, 10).Value = 292  # Robbery: 290 -> 292
$ws.Cells.Item(7, 10).Value = 1010  # Total: 1004 -> 1010

# --- Little Village ---
$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(3, 10).Value = 180  # Aggravated Battery: 179 -> 180
$ws.Cells.Item(7, 10).Value = 613  # Total: 612 -> 613
